# The target diff reorders two pairs of data rows in the "Artfynd" sheet:
#   - row 19 and row 20 trade places (entire record, every column)
#   - row 21 and row 22 trade places (entire record, every column)
# Nothing else in the sheet changes. We implement this as two in-place row
# swaps using whole-row range Value() reads/writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies columns A..AY (51 columns).
$lastCol = 51

# The Startdatum/Slutdatum columns (Y and AA) store dates as plain text
# (e.g. "2023-08-18") rather than real Excel dates. Force those columns to
# Text format first so that re-assigning the read-back values does not let
# Excel auto-convert the strings into date serial numbers.
$ws.Range("Y19:Y22").NumberFormat = "@"
$ws.Range("AA19:AA22").NumberFormat = "@"

# --- Swap row 19 <-> row 20 ---
$row19 = $ws.Range($ws.Cells.Item(19, 1), $ws.Cells.Item(19, $lastCol))
$row20 = $ws.Range($ws.Cells.Item(20, 1), $ws.Cells.Item(20, $lastCol))

$values19 = $row19.Value()
$values20 = $row20.Value()

$row19.Value = $values20
$row20.Value = $values19

# --- Swap row 21 <-> row 22 ---
$row21 = $ws.Range($ws.Cells.Item(21, 1), $ws.Cells.Item(21, $lastCol))
$row22 = $ws.Range($ws.Cells.Item(22, 1), $ws.Cells.Item(22, $lastCol))

$values21 = $row21.Value()
$values22 = $row22.Value()

$row21.Value = $values22
$row22.Value = $values21
